$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing timestamp values (rows 2, 3, 5) to the new run's timestamps.
# Row 4 (Wed Nov 01 15:53:36 EDT 2023 / DONOTRUN) is left untouched.
$ws.Range("B2").Value = "Thu Jan 25 17:47:19 EST 2024"
$ws.Range("B3").Value = "Thu Jan 25 17:47:31 EST 2024"
$ws.Range("B5").Value = "Thu Jan 25 17:47:43 EST 2024"

# Append two new rows of test data for the "Estate Tax" tax type.
$ws.Range("D6").Value = "Existing Liability w/Notice Number"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Borders.LineStyle = 1
$ws.Range("D6").WrapText = $true

$ws.Range("E6").Value = "Estate Tax"

$ws.Range("D7").Value = "New Tax Return Amount Due"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Borders.LineStyle = 1
$ws.Range("D7").WrapText = $true

$ws.Range("E7").Value = "Estate Tax"

# Move the active selection to the last entered cell.
$ws.Range("E7").Select()
